# Update the lot/date tracking info and received-quantity figures across the
# four PR_GAST_3 stock-card sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: PR_GAST_3_Uncut_Sheet -----------------------------------------
$wsUncut = $wb.Worksheets.Item("PR_GAST_3_Uncut_Sheet")

# "08.08.2022" is ambiguous as day/month (both 08), so Excel would otherwise
# auto-convert it to a date serial; the leading apostrophe keeps it literal
# text, same as a user typing it in the formula bar.
$wsUncut.Range("A12").Value = "'08.08.2022"
$wsUncut.Range("B12").Value = 20
$wsUncut.Range("D12").Value = 20
$wsUncut.Range("G12").Value = "2024-06"

# --- Sheet 2: PR_GAST_3_Cassette --------------------------------------------
$wsCassette = $wb.Worksheets.Item("PR_GAST_3_Cassette")

$wsCassette.Range("A12").Value = "'08.08.2022"
$wsCassette.Range("B12").Value = 1314
$wsCassette.Range("C12").Value = 1300
$wsCassette.Range("E12").Value = 1315
$wsCassette.Range("F12").Value = 1300

# Move the cursor to H12 on this sheet (it was G12 before).
$wsCassette.Activate()
$wsCassette.Range("H12").Select()

# --- Sheet 3: PR_GAST_3_Pipette ---------------------------------------------
$wsPipette = $wb.Worksheets.Item("PR_GAST_3_Pipette")

$wsPipette.Range("B12").Value = 4162
$wsPipette.Range("D12").Value = 1300

# --- Restore the originally active sheet/tab --------------------------------
$wsBuffer = $wb.Worksheets.Item("PR_GAST_3_Buffer")
$wsBuffer.Activate()
